$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Fuente: XIII Censo de Población y Vivienda 2010, Encuesta
# Intercensal 2015 y Censo de Población y Vivienda 2020, INEGI." footnote
# that lived in the merged cell A7:G7 (leaves the cell formatting intact,
# only drops the text, same as the source edit).
$ws.Range("A7:G7").Value = $null

# Move / restore the active selection to D15 (previously C17).
$ws.Range("D15").Select()
